$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(15, 1000, 10, 0.2, "'True"),
    @(16, 1000, 10, 0.2, "'True"),
    @(18, 1000, 10, 0.2, "'True"),
    @(19, 1000, 10, 0.2, "'True"),
    @(20, 1000, 10, 0.2, "'True"),
    @(21, 1000, 10, 0.2, "'True"),
    @(17, 10, 10, 0.2, "'True")
)

$row = 3
foreach ($rowData in $data) {
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $row++
}

$ws.Range("B9").Select()
